$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.75
$ws.Range("B3").Value = 0.07000000000000001
$ws.Range("B4").Value = 0.73
$ws.Range("B5").Value = 0.73
$ws.Range("B6").Value = 0.6
$ws.Range("B7").Value = -0.38
$ws.Range("B8").Value = 0.32
$ws.Range("B9").Value = -0.23
$ws.Range("B10").Value = -0.46
$ws.Range("B11").Value = 0.2
